$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column (C) for rows 2-6 from 2023-10-25 to 2023-11-03
$newDate = Get-Date -Year 2023 -Month 11 -Day 3 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
